$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Column Type")
$ws.Activate()
